$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    12 = -21.626
    27 = -21.718
    32 = -21.667
    36 = -20.339
    38 = -19.741
    46 = -21.858
    54 = -21.862
    55 = -22.247
    56 = -22.095
    67 = -21.495
    69 = -21.721
    72 = -21.445
    83 = -21.996
    86 = -22.093
    91 = -21.522
    93 = -21.259
    99 = -20.43
}

foreach ($row in $updates.Keys) {
    $ws.Range("A$row").Value = $updates[$row]
}
